$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Brand"
$ws.Range("B1").Value = "Customer"

# Rows 2-10: group "1"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = "brandfile1"
    $ws.Cells.Item($r, 2).Value = "customerfiel1"
}

# Rows 11-19: group "2"
for ($r = 11; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = "brandfile2"
    $ws.Cells.Item($r, 2).Value = "customerfiel2"
}
